$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New consolidated trigger array replaces the two different trigger arrays
# that used to be split between rows 8-19 and rows 20-21.
$newTrigger = '["01", "04", "08", "13", "17", "22", "26", "30"]'
$ws.Range("C8:C21").Value = $newTrigger

# Rows 20/21 used to carry the 18:15/18:20 time slots together with the
# old "Мы, Единый Народ России..." trigger text; they are replaced with a
# new pair of 11:55/12:00 time slots and new wording for the appeal text.
$ws.Range("A20").Value = "Мы, Единый Народ России, обращаемся, к Создателю Мира сего, с просьбой, поручить нам управление Планетой в Шестой Эпохе, через Белого Царя, ибо формируется, Навечно, сотворчество людей-Богов с Создателем!"
$ws.Range("B20").Value = "11:55 - 11:59"
$ws.Range("B21").Value = "12:00 - 12:04"

# The per-row numbering in column E is no longer used.
$ws.Range("E8:E21").ClearContents()

# Row 20 no longer needs the extra-tall wrapped height of the old, much
# longer paragraph.
$ws.Rows(20).RowHeight = 45

# Update the visible selection/scroll position left behind by the editor.
$ws.Range("A24:A26").Select()
